$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Shape "Content Placeholder 2": split the existing sentence into two
# runs, then append a brand-new paragraph (made up of several runs) about
# lightsaber duels / sightlabers.
$content = $s.Shapes.Item(2)
$tr = $content.TextFrame.TextRange

$tr.Text = "The game will be built using the Unity game engine, and run on any VR-capable "
$tr.InsertAfter("computer") | Out-Null
$tr.InsertAfter("`rCome on, I know you had lightsaber duels as a kid. ") | Out-Null
$tr.InsertAfter("Erm") | Out-Null
$tr.InsertAfter(", I mean, ") | Out-Null
$tr.InsertAfter("sightlabers") | Out-Null
$tr.InsertAfter(". Yep. It’s just that, but in VR.") | Out-Null

# --- Shape "TextBox 3": widen the box to fit the new caption text and swap
# the text itself.
$caption = $s.Shapes.Item(3)
$caption.Width = 6401689 / 12700
$caption.TextFrame.TextRange.Text = "You know, in blindfold mode, I could just play sounds and “fake it”."
